# Stimuli.xlsx - "Changed stimuli pattern for better stability"
# Apply the new stimuli input values on Sheet1; every other changed cell in
# the diff (row 1 series, rows 28/29/30/32/35/37/39) is a formula that
# recalculates automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B1: -3 -> -4 (drives the whole shared-formula series in row 1)
$ws.Range("B1").Value = -4

# Row 15 (binary switch row)
$ws.Range("O15").Value = 2

# Row 17
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0

# Row 18
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0

# Row 23
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0

# Move the active selection from X17 to Q11, matching the saved view state
$ws.Activate()
$ws.Range("Q11").Select()
